# rename hands-on to labs
$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 18: Title "Hands On Exercise:\nCreate an EC2 Instance"
#           -> "Lab: Create an EC2 Instance" (single line, no break)
#           and drop the normAutofit fontScale override.
# ---------------------------------------------------------------------
$s18 = $p.Slides.Item(18)
$titleShp = $s18.Shapes.Item(2)

# Reset autofit so PowerPoint stops forcing a 90% font scale
# (<a:normAutofit fontScale="90000"/> -> <a:normAutofit/>)
$titleShp.TextFrame.AutoSize = 2

$titleRange = $titleShp.TextFrame.TextRange

# Replace the whole range (this also removes the <a:br/> soft line break)
$wholeTitle = $titleRange.Characters(1, $titleRange.Text.Length)
$wholeTitle.Text = "Lab: Create an EC2 Instance"

# Re-select fresh range and split into the runs matching the authored edit
$titleRange2 = $titleShp.TextFrame.TextRange
$titleRange2.Characters(1, 4).Text = "Lab:"
$titleRange2.Characters(5, 1).Text = " "
$titleRange2.Characters(6, 7).Text = "Create "
$titleRange2.Characters(13, 15).Text = "an EC2 Instance"

# ---------------------------------------------------------------------
# Slide 8: Content placeholder, 2nd paragraph (lvl 1)
#          "11 " + "Available Regions (as of " + "November, " + "2015" + ")"
#          -> single run "11 Available Regions (as of November, 2015)"
# ---------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$contentShp = $s8.Shapes.Item(1)
$contentRange = $contentShp.TextFrame.TextRange

$fullText = $contentRange.Text
$target = "11 Available Regions (as of November, 2015)"
$startIdx = $fullText.IndexOf($target)
$regionsRange = $contentRange.Characters($startIdx + 1, $target.Length)
$regionsRange.Text = $target
